$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value2 = 'SELECT
    COUNT(DISTINCT p.program_acronym) AS "Programs",
    COUNT(DISTINCT st.clinical_study_designation) AS "Studies",
    COUNT(DISTINCT c.case_record_id) AS "Cases",
    COUNT(DISTINCT smp.sample_id) AS "Samples",
    COUNT(DISTINCT cf.file_name) AS "Case Files",         
    COUNT(DISTINCT sf.file_name) AS "Study Files" 
FROM 
    df_program p
JOIN 
    df_study st ON p.program_acronym = st."program.program_acronym"
JOIN 
    df_case c ON st.clinical_study_designation = c."study.clinical_study_designation"
JOIN 
    df_demographic dmg ON dmg."case.case_record_id" = c.case_record_id
JOIN 
    df_diagnosis diag ON diag."case.case_record_id" = c.case_record_id
JOIN 
    df_enrollment enr ON enr."case.case_record_id" = c.case_record_id
JOIN 
    df_sample smp ON smp."case.case_record_id" = c.case_record_id
LEFT JOIN 
    df_case_file cf ON cf."sample.sample_id" = smp.sample_id
LEFT JOIN 
    df_study_file sf ON sf."study.clinical_study_designation" = st.clinical_study_designation
WHERE 
    st.clinical_study_designation = ''MGT01'' AND dmg.neutered_indicator = ''Yes'';'

$ws.Range("B2").Value2 = 'SELECT
    DISTINCT c.case_record_id AS "Case ID",
    st.clinical_study_designation AS "Study Code",
    st.clinical_study_type AS "Study Type",
    dmg.breed AS "Breed",
    diag.disease_term AS "Diagnosis",
    diag.stage_of_disease AS "Stage Of Disease",
    CASE 
    WHEN dmg.patient_age_at_enrollment = CAST(dmg.patient_age_at_enrollment AS INT) 
    THEN CAST(CAST(dmg.patient_age_at_enrollment AS INT) AS TEXT)
    ELSE CAST(dmg.patient_age_at_enrollment AS TEXT)
END AS "Age", 
    COALESCE(TRIM(dmg.sex), '''') AS "Sex",   
    COALESCE(dmg.neutered_indicator, '''') AS "Neutered Status",  
    COALESCE(
        CASE 
            WHEN dmg.weight = CAST(dmg.weight AS INT) THEN CAST(dmg.weight AS INT)
            ELSE dmg.weight 
        END, 
    '''') AS "Weight (kg)",        
    COALESCE(diag.best_response, '''') AS "Response to Treatment",  
    COALESCE(coh.cohort_description, '''') AS "Cohort"  
FROM 
    df_case c
JOIN 
    df_study st ON c."study.clinical_study_designation" = st.clinical_study_designation
JOIN 
    df_demographic dmg ON dmg."case.case_record_id" = c.case_record_id
JOIN 
    df_diagnosis diag ON diag."case.case_record_id" = c.case_record_id
JOIN 
    df_enrollment enr ON enr."case.case_record_id" = c.case_record_id
JOIN 
    df_program p ON st."program.program_acronym" = p.program_acronym
JOIN 
    df_sample smp ON smp."case.case_record_id" = c.case_record_id
JOIN 
    df_publication pub ON pub."study.clinical_study_designation" = st.clinical_study_designation
LEFT JOIN 
    df_case_file cf ON cf."sample.sample_id" = smp.sample_id
LEFT JOIN 
    df_study_file sf ON sf."study.clinical_study_designation" = st.clinical_study_designation
LEFT JOIN
    df_cohort coh ON coh."study.clinical_study_designation" = st.clinical_study_designation
WHERE
    st.clinical_study_designation = ''MGT01'' AND dmg.neutered_indicator = ''Yes''
ORDER BY 
    c.case_record_id ASC
LIMIT 100;'

$ws.Range("B3").Value2 = 'SELECT DISTINCT
    smp.sample_id AS "Sample ID",
    c.case_record_id AS "Case ID",
    COALESCE(dmg.breed, '''') AS "Breed",
    COALESCE(diag.disease_term, '''') AS "Diagnosis",
    COALESCE(smp.sample_site, '''') AS "Sample Site",
    COALESCE(smp.summarized_sample_type, '''') AS "Sample Type",
    COALESCE(smp.specific_sample_pathology, '''') AS "Pathology/Morphology",
    COALESCE(smp.tumor_grade, '''') AS "Tumor Grade",
    COALESCE(smp.sample_chronology, '''') AS "Sample Chronology",
    COALESCE(smp.percentage_tumor, '''') AS "Percentage Tumor",
    COALESCE(smp.necropsy_sample, '''') AS "Necropsy Sample",
    COALESCE(smp.sample_preservation, '''') AS "Sample Preservation"
FROM 
    df_sample smp
JOIN 
    df_case c ON smp."case.case_record_id" = c.case_record_id
JOIN 
    df_publication pub ON pub."study.clinical_study_designation" = st.clinical_study_designation
JOIN 
    df_demographic dmg ON dmg."case.case_record_id" = c.case_record_id
JOIN 
    df_diagnosis diag ON diag."case.case_record_id" = c.case_record_id
JOIN 
    df_enrollment enr ON enr."case.case_record_id" = c.case_record_id
JOIN 
    df_program p ON st."program.program_acronym" = p.program_acronym
JOIN 
    df_study st ON c."study.clinical_study_designation" = st.clinical_study_designation
LEFT JOIN 
    df_case_file cf ON cf."sample.sample_id" = smp.sample_id
LEFT JOIN 
    df_study_file sf ON sf."study.clinical_study_designation" = st.clinical_study_designation
WHERE 
   st.clinical_study_designation = ''MGT01'' AND dmg.neutered_indicator = ''Yes''
ORDER BY 
    smp.sample_id ASC
LIMIT 100;'

$ws.Range("B4").Value2 = 'SELECT 
    DISTINCT cf.file_name AS "File Name",
    CASE
        WHEN cf.file_name LIKE ''%.bai'' THEN ''bai''
        WHEN cf.file_name LIKE ''%.bam'' THEN ''bam''
        WHEN cf.file_name LIKE ''%.csv'' THEN ''csv''
        WHEN cf.file_name LIKE ''%.doc'' THEN ''doc''
        WHEN cf.file_name LIKE ''%.docx'' THEN ''docx''
        WHEN cf.file_name LIKE ''%.gz'' THEN ''gz''
        WHEN cf.file_name LIKE ''%.pdf'' THEN ''pdf''
        WHEN cf.file_name LIKE ''%.rtf'' THEN ''rtf''
        WHEN cf.file_name LIKE ''%.tbi'' THEN ''tbi''
        WHEN cf.file_name LIKE ''%.tif'' THEN ''tif''
        WHEN cf.file_name LIKE ''%.xls'' THEN ''xls''
        WHEN cf.file_name LIKE ''%.xlsx'' THEN ''xlsx''
        ELSE ''Unknown''
    END AS "Format",
    cf.file_type AS "File Type",
    CASE     
    WHEN cf.file_size >= 1024 * 1024 * 1024 THEN 
        CASE 
            WHEN ROUND(cf.file_size / (1024.0 * 1024.0 * 1024.0), 2) = CAST(ROUND(cf.file_size / (1024.0 * 1024.0 * 1024.0), 0) AS INT) 
            THEN CAST(CAST(ROUND(cf.file_size / (1024.0 * 1024.0 * 1024.0), 0) AS INT) AS TEXT) || '' GB''
            ELSE ROUND(cf.file_size / (1024.0 * 1024.0 * 1024.0), 2) || '' GB''
        END
    WHEN cf.file_size >= 1024 * 1024 THEN 
        CASE 
            WHEN ROUND(cf.file_size / (1024.0 * 1024.0), 2) = CAST(ROUND(cf.file_size / (1024.0 * 1024.0), 0) AS INT) 
            THEN CAST(CAST(ROUND(cf.file_size / (1024.0 * 1024.0), 0) AS INT) AS TEXT) || '' MB''
            ELSE ROUND(cf.file_size / (1024.0 * 1024.0), 2) || '' MB''
        END
    WHEN cf.file_size >= 1024 THEN 
        CASE 
            WHEN ROUND(cf.file_size / 1024.0, 2) = CAST(ROUND(cf.file_size / 1024.0, 0) AS INT) 
            THEN CAST(CAST(ROUND(cf.file_size / 1024.0, 0) AS INT) AS TEXT) || '' KB''
            ELSE ROUND(cf.file_size / 1024.0, 2) || '' KB''
        END
    ELSE 
        CASE 
            WHEN ROUND(cf.file_size, 2) = CAST(ROUND(cf.file_size, 0) AS INT) 
            THEN CAST(CAST(ROUND(cf.file_size, 0) AS INT) AS TEXT) || '' Bytes''
            ELSE ROUND(cf.file_size, 2) || '' Bytes''
        END
END AS "Size",
    ''sample'' AS "Association",
    cf.file_description AS "Description",
    smp.sample_id AS "Sample ID",
    c.case_record_id AS "Case ID",
    dmg.breed AS "Breed",
    diag.disease_term AS "Diagnosis"
FROM 
    df_case_file cf
JOIN 
    df_sample smp ON cf."sample.sample_id" = smp.sample_id
JOIN 
    df_case c ON smp."case.case_record_id" = c.case_record_id
JOIN 
    df_study st ON c."study.clinical_study_designation" = st.clinical_study_designation
JOIN 
    df_program p ON st."program.program_acronym" = p.program_acronym
JOIN 
    df_demographic dmg ON dmg."case.case_record_id" = c.case_record_id
JOIN 
    df_diagnosis diag ON diag."case.case_record_id" = c.case_record_id
JOIN 
    df_enrollment enr ON enr."case.case_record_id" = c.case_record_id
JOIN 
    df_publication pub ON pub."study.clinical_study_designation" = st.clinical_study_designation
LEFT JOIN 
    df_study_file sf ON sf."study.clinical_study_designation" = st.clinical_study_designation
WHERE
    st.clinical_study_designation = ''MGT01'' AND dmg.neutered_indicator = ''Yes''
ORDER BY 
    cf.file_name ASC
LIMIT 100;'

$ws.Range("B5").Value2 = 'SELECT DISTINCT
    sf.file_name AS "File Name",
    sf.file_type AS "File Type",
    ''study'' AS "Association",
    sf.file_description AS "Description",
    CASE
        WHEN sf.file_name LIKE ''%.bai'' THEN ''bai''
        WHEN sf.file_name LIKE ''%.bam'' THEN ''bam''
        WHEN sf.file_name LIKE ''%.csv'' THEN ''csv''
        WHEN sf.file_name LIKE ''%.doc'' THEN ''doc''
        WHEN sf.file_name LIKE ''%.docx'' THEN ''docx''
        WHEN sf.file_name LIKE ''%.gz'' THEN ''gz''
        WHEN sf.file_name LIKE ''%.pdf'' THEN ''pdf''
        WHEN sf.file_name LIKE ''%.rtf'' THEN ''rtf''
        WHEN sf.file_name LIKE ''%.tbi'' THEN ''tbi''
        WHEN sf.file_name LIKE ''%.tif'' THEN ''tif''
        WHEN sf.file_name LIKE ''%.xls'' THEN ''xls''
        WHEN sf.file_name LIKE ''%.xlsx'' THEN ''xlsx''
        ELSE ''Unknown''
    END AS "Format",
    CASE     
        WHEN sf.file_size >= 1024 * 1024 * 1024 THEN 
            ROUND(sf.file_size / (1024.0 * 1024.0 * 1024.0), 2) || '' GB'' 
        WHEN sf.file_size >= 1024 * 1024 THEN 
            ROUND(sf.file_size / (1024.0 * 1024.0), 2) || '' MB'' 
        WHEN sf.file_size >= 1024 THEN 
            ROUND(sf.file_size / 1024.0, 2) || '' KB'' 
        ELSE 
            ROUND(sf.file_size, 2) || '' Bytes'' 
    END AS "Size",
    st.clinical_study_designation AS "Study Code"
FROM 
    df_case_file cf
JOIN 
    df_sample smp ON cf."sample.sample_id" = smp.sample_id
JOIN 
    df_case c ON smp."case.case_record_id" = c.case_record_id
JOIN 
    df_study st ON c."study.clinical_study_designation" = st.clinical_study_designation
JOIN 
    df_program p ON st."program.program_acronym" = p.program_acronym
JOIN 
    df_demographic dmg ON dmg."case.case_record_id" = c.case_record_id
JOIN 
    df_diagnosis diag ON diag."case.case_record_id" = c.case_record_id
JOIN 
    df_enrollment enr ON enr."case.case_record_id" = c.case_record_id
JOIN 
    df_publication pub ON pub."study.clinical_study_designation" = st.clinical_study_designation
LEFT JOIN 
    df_study_file sf ON sf."study.clinical_study_designation" = st.clinical_study_designation
WHERE
    st.clinical_study_designation = ''MGT01'' AND dmg.neutered_indicator = ''Yes''
ORDER BY 
    sf.file_name ASC
LIMIT 100;'

$ws.Range("C5").Select()
